$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.811.42"
$ws.Range("E2").Value = "  +2.97%  "

$ws.Range("D3").Value = "2.267.22"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.39"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.66"
$ws.Range("E6").Value = "  +5.13%  "

$ws.Range("E7").Value = "  -1.50%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.509"
$ws.Range("E9").Value = "  -1.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.21"
$ws.Range("E10").Value = "  +2.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.08"
$ws.Range("E12").Value = "  -2.86%  "

$ws.Range("E13").Value = "  -1.22%  "

$ws.Range("D14").Value = "2.612.86"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").Value = "2.266.98"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.63"
$ws.Range("E16").Value = "  -0.70%  "

$ws.Range("D17").Value = "46.791.40"
$ws.Range("E17").Value = "  +3.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.793"
$ws.Range("E18").Value = "  -2.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.72"
$ws.Range("E19").Value = "  -3.65%  "

$ws.Range("D20").Value = "0.0₃0964"
$ws.Range("E20").Value = "  +4.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.81"
$ws.Range("E21").Value = "  -4.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.59"
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.66"
$ws.Range("E23").Value = "  +3.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.80"
$ws.Range("E24").Value = "  -3.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.86"
$ws.Range("E26").Value = "  -3.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.67"
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("E28").Value = "  -1.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.62"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.22"
$ws.Range("E30").Value = "  +2.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.81"
$ws.Range("E31").Value = "  +8.69%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "146.38"
$ws.Range("E32").Value = "  -3.34%  "

$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.31"
$ws.Range("E33").Value = "  +11.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.37"
$ws.Range("E34").Value = "  -3.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0769"
$ws.Range("E35").Value = "  -3.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").Value = "  +9.42%  "

$ws.Range("E37").Value = "  -2.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.76"
$ws.Range("E38").Value = "  +15.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.69"
$ws.Range("E39").Value = "  -5.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.85"
$ws.Range("E40").Value = "  -3.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0298"
$ws.Range("E41").Value = "  -5.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.11"
$ws.Range("E42").Value = "  -4.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.90"
$ws.Range("E44").Value = "  -2.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.02"
$ws.Range("E45").Value = "  +18.15%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.781.44"
$ws.Range("E46").Value = "  +0.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "71.27"
$ws.Range("E47").Value = "  +0.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.185"
$ws.Range("E48").Value = "  -4.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.81"
$ws.Range("E49").Value = "  +0.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "94.15"
$ws.Range("E50").Value = "  -2.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.82"
$ws.Range("E51").Value = "  -1.60%  "
